$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row 2 (shifting the current rows 2-5 down to 3-6)
$ws.Rows.Item(2).Insert()

# Fill the new "slug" row with the short local identifiers for each column
$ws.Range("A2").Value = "poca-limpieza"
$ws.Range("B2").Value = "pocas-zonas-verdes"
$ws.Range("C2").Value = "malas-comunicaciones"
$ws.Range("D2").Value = "delincuencia-zona"
$ws.Range("E2").Value = "numero-viviendas"
$ws.Range("F2").Value = "contaminacion"
$ws.Range("G2").Value = "falta-de-servicios-de-aseo"
$ws.Range("H2").Value = "aragon"
$ws.Range("I2").Value = "ruidos-exteriores"
